$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Acervo Bibliográfico")

# Fix the header typo: "Edicao" -> "Edição"
$ws.Range("J1").Value = "Edição"

# Update the active selection to match the new state
$ws.Range("J2").Select()
